# Edit scenariolibrary.xlsx - update BF scenario rows
# "numbers suck but trying to fix them"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 used to hold the lone "china-HeEtAl2017" reference row (A: source tag,
# B: coal, C: air, D: BF+BOF). It is repurposed into the first of a new block
# of EU/China blast-furnace scenario name stubs, and the now-stale B/C/D
# values in that row are removed entirely.
$ws.Range("A34").Value = "EU-BF-I"
$ws.Range("B34:D34").ClearFormats()
$ws.Range("B34:D34").ClearContents()

# New rows 35-41: additional scenario name stubs (column A only).
$ws.Range("A35").Value = "EU-BF-C"
$ws.Range("A36").Value = "EU-BF-M"
$ws.Range("A37").Value = "EU-BF-F"
$ws.Range("A38").Value = "China-BF-I"
$ws.Range("A39").Value = "China-BF-C"
$ws.Range("A40").Value = "China-BF-M"
$ws.Range("A41").Value = "China-BF-F"
